$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()
$ws.Range("A1:E1").Borders.Item(9).LineStyle = 1
$ws.Range("A1:E1").Borders.Item(9).Weight = 4
